# Update scraped crypto market data (price / 1h volume change) and fix
# two swapped rows (FTXToken / GateToken) to match the refreshed source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.13"
$ws.Range("E2").Value = "'-5.25%"
$ws.Range("E3").Value = "'-2.81%"
$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'-4.00%"
$ws.Range("D5").Value = "'0.07468"
$ws.Range("E5").Value = "'-2.67%"
$ws.Range("D6").Value = "'7.731"
$ws.Range("E6").Value = "'-1.84%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.682"
$ws.Range("E7").Value = "'4.95%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.783"
$ws.Range("E8").Value = "'2.00%"
$ws.Range("D9").Value = "'0.9292"
$ws.Range("E9").Value = "'1.13%"
$ws.Range("D10").Value = "'0.1691"
$ws.Range("E10").Value = "'-2.44%"
$ws.Range("D11").Value = "'0.07167"
$ws.Range("E11").Value = "'-5.38%"
$ws.Range("D12").Value = "'0.07975"
$ws.Range("E12").Value = "'-3.67%"
$ws.Range("D13").Value = "'0.03031"
$ws.Range("E13").Value = "'0.10%"
$ws.Range("D14").Value = "'0.09905"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'-1.52%"
$ws.Range("D16").Value = "'0.006311"
$ws.Range("E16").Value = "'0.69%"
$ws.Range("D17").Value = "'3.456"
$ws.Range("E17").Value = "'-0.68%"
$ws.Range("D18").Value = "'2.226"
$ws.Range("E18").Value = "'-0.61%"
$ws.Range("E19").Value = "'-0.99%"
$ws.Range("D20").Value = "'0.1350"
$ws.Range("E20").Value = "'2.37%"
$ws.Range("D21").Value = "'4.569"
$ws.Range("E21").Value = "'7.53%"
$ws.Range("D22").Value = "'0.04644"
$ws.Range("E22").Value = "'1.70%"
$ws.Range("D23").Value = "'0.1555"
$ws.Range("E23").Value = "'-4.39%"
$ws.Range("E24").Value = "'-0.15%"
$ws.Range("D25").Value = "'0.004420"
$ws.Range("E25").Value = "'-1.80%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.41%"
$ws.Range("D27").Value = "'0.0001879"
$ws.Range("E27").Value = "'7.99%"
$ws.Range("D39").Value = "'0.01661"
$ws.Range("E39").Value = "'-3.09%"
$ws.Range("D40").Value = "'0.04441"
$ws.Range("E40").Value = "'-4.34%"
$ws.Range("D41").Value = "'0.007049"
$ws.Range("E41").Value = "'-2.27%"
$ws.Range("D42").Value = "'0.1325"
$ws.Range("E42").Value = "'-3.28%"
$ws.Range("D43").Value = "'0.002074"
$ws.Range("E43").Value = "'-8.03%"
$ws.Range("D44").Value = "'0.01233"
$ws.Range("E44").Value = "'-15.32%"
$ws.Range("D45").Value = "'0.00005996"
$ws.Range("E45").Value = "'-3.07%"
$ws.Range("D46").Value = "'1.918"
$ws.Range("E46").Value = "'1.34%"
$ws.Range("D47").Value = "'0.01102"
$ws.Range("E47").Value = "'-15.18%"
